$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.075935482978821
$ws.Range("B1").Value = 0.9521679878234863
$ws.Range("C1").Value = 6.583098411560059
$ws.Range("D1").Value = 2.019391775131226
$ws.Range("E1").Value = 1.122024297714233
